$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.89%"

# Row 3
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "8.82%"

# Row 4
$ws.Range("B4").Value = "LEO"
$ws.Range("C4").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "3.597"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.32%"

# Row 5
$ws.Range("B5").Value = "HuobiToken"
$ws.Range("C5").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "5.098"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.11%"

# Row 6
$ws.Range("B6").Value = "Cronos"
$ws.Range("C6").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.08140"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.48%"

# Row 7
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.961"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.74%"

# Row 8
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.200"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.01%"

# Row 9
$ws.Range("B9").Value = "KuCoinToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.934"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.98%"

# Row 10
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9280"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.76%"

# Row 11
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1437"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "11.68%"

# Row 12
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1951"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.96%"

# Row 13
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09109"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.08%"

# Row 14
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03506"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.54%"

# Row 15
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09834"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.19%"

# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001412"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.41%"

# Row 17
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006067"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.72%"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.475"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.12%"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3447"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.06%"

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.07%"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.804"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-7.04%"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2434"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-6.51%"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04457"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.87%"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001239"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.38%"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004856"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.98%"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001302"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.03%"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02103"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "7.76%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05098"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-8.05%"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007467"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.23%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01014"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.60%"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1363"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.66%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002143"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.44%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01041"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.89%"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006217"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.12%"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.04%"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003063"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001603"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.46%"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.04%"
